# Update Master Data Tables - machine_master sheet:
#  - Reformat existing mac_address values (rows 2-21) from lower-case
#    colon-separated to upper-case dash-separated form.
#  - Append 9 new machine rows (22-30 / id 10021-10029, "Machine 21".."Machine 29").
#  - Adjust column C width and the saved selection to match the new sheet extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  "Machine 1","Machine 2","Machine 3","Machine 4","Machine 5","Machine 6","Machine 7","Machine 8","Machine 9","Machine 10",
  "Machine 11","Machine 12","Machine 13","Machine 14","Machine 15","Machine 16","Machine 17","Machine 18","Machine 19","Machine 20",
  "Machine 21","Machine 22","Machine 23","Machine 24","Machine 25","Machine 26","Machine 27","Machine 28","Machine 29"
)

$macs = @(
  "8C-16-45-5A-5D-0D","8C-16-45-88-E1-0D","00-FF-D3-E3-9A-27","8C-16-45-5A-62-41","E8-6A-64-1D-75-E4",
  "8C-16-45-FA-94-B7","8C-16-45-1A-0F-62","E8-6A-64-1C-52-6E","48-51-B7-10-35-A6","8C-16-45-38-F3-F3",
  "D4-3D-7E-58-CC-45","8C-16-45-5A-5D-96","8C-16-45-5A-5D-8E","8C-16-45-33-A5-5F","3C-95-09-F9-EA-DF",
  "8C-16-45-88-E7-0B","B4-69-21-5A-DB-C4","E8-6A-64-1D-48-B7","8C-16-45-59-69-09 ","98-E7-F4-30-16-5A ",
  "38-BA-F8-53-C7-8F","E8-6A-64-1C-58-C2","E4-A4-71-CE-BA-93","54-E1-AD-EA-30-C9","8C-16-45-65-DD-40",
  "58-20-B1-D6-C3-BE","8C-16-45-38-F0-25","6C-88-14-AC-EF-55","3C-6A-A7-C0-DF-27"
)

$serials = @(
  "NM5328114630","WY2132605316","CM6384145127","NT894252578","YM866672706","WT6501645780",
  "LK8186452621","NR3264783870","RW437027336","SI158158531","XF3416823469","BW4524978011",
  "DB289579153","SI4597903231","TJ7809002958","JR6082789079","SA3722889241","RR2683722548",
  "PO6528391346","FB5962911652","FB5962911653","FB5962911654","FB5962911655","FB5962911656",
  "FB5962911657","FB5962911658","FB5962911659","FB5962911661","FB5962911662"
)

$ips = @(
  "192.168.0.150","192.168.0.133","192.168.0.161","192.168.0.259","192.168.0.119","192.168.0.177",
  "192.168.0.227","192.168.0.207","192.168.0.220","192.168.0.242","192.168.0.173","192.168.0.203",
  "192.168.0.112","192.168.0.178","192.168.0.267","192.168.0.149","192.168.0.127","192.168.0.248",
  "192.168.0.121","192.168.0.215","192.168.0.874","192.168.0.721","192.168.0.841","192.168.0.186",
  "192.168.0.627","192.168.0.879","192.168.0.628","192.168.0.306","192.168.0.355"
)

for ($i = 0; $i -lt 29; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = 10001 + $i
  $ws.Cells.Item($row, 2).Value = $names[$i]
  $ws.Cells.Item($row, 3).Value = $macs[$i]
  $ws.Cells.Item($row, 4).Value = $serials[$i]
  $ws.Cells.Item($row, 5).Value = $ips[$i]
  $ws.Cells.Item($row, 6).Value = 1001
  $ws.Cells.Item($row, 7).Value = "eng"
  $ws.Cells.Item($row, 8).Value = $true
  $ws.Cells.Item($row, 9).Value = "superadmin"
  $ws.Cells.Item($row, 10).Value = "now()"
}

# Column C (mac_address) is now a touch wider to fit the new dash-separated values.
$ws.Columns.Item(3).ColumnWidth = 16.15

# Leave the saved selection where Excel would land after entering data through row 30.
$ws.Range("A31:XFD1048576").Select() | Out-Null
